# Insert a new "l1" column right before the existing "nota_iniciativa" column
# (which currently lives in column Q). This shifts "nota_iniciativa" one
# column to the right (column R) and frees up column Q for the new "l1"
# data, matching the pattern used by the neighbouring "c1" column (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift column Q ("nota_iniciativa") and everything after it one column to
# the right, inserting a brand-new, empty column Q.
$ws.Columns("Q").Insert()

# Header for the newly inserted column.
$ws.Range("Q1").Value = "l1"

# The source data rows run from row 2 through row 85; fill them with the
# default value 0, same as the neighboring "c1" (P) and "nota_iniciativa"
# (now R) columns.
$ws.Range("Q2:Q85").Value = 0
